$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the two employee records in rows 16 and 17, keeping formatting intact.
# Row 16 currently holds GABRIEL EDUARDO BAUTISTA MANGA's data; row 17 holds
# ORLANDO MANUEL QUIÑONES PAYARES's data. The updated account-statement
# database reorders them so Orlando's record comes first (row 16) and
# Gabriel's record comes second (row 17).

$ws.Range("C16").Value = "1047475480"
$ws.Range("D16").Value = "ORLANDO MANUEL QUIÑONES PAYARES"
$ws.Range("E16").Value = "2501"
$ws.Range("F16").Value = 2800
$ws.Range("G16").Value = 2100000

$ws.Range("C17").Value = "1047458118"
$ws.Range("D17").Value = "GABRIEL EDUARDO BAUTISTA MANGA"
$ws.Range("E17").Value = "2502"
$ws.Range("F17").Value = 2667
$ws.Range("G17").Value = 2000000
